# Fix processing of borderline number of pictures per card.
# Some products had fewer than 5 images; trailing placeholder/duplicate
# image rows for those products need to be removed so each product's
# image rows match its real image count.
#
# snowboard102 -> keep only 1 image row (remove rows 13-16)
# snowboard105 -> keep only 3 image rows (remove rows 30-31)
# snowboard108 -> keep only 3 image rows (remove rows 45-46)
#
# Delete from the bottom of the sheet upward so earlier row numbers stay
# valid while later ranges are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45:H46").EntireRow.Delete()
$ws.Range("A30:H31").EntireRow.Delete()
$ws.Range("A13:H16").EntireRow.Delete()

# Restore the explicit column widths for column A (key) and column H
# (variants.images.dimensions.h) that ship with the fixed workbook.
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 5
